$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.627.38'
$ws.Range('E2').Value = '  +1.49%  '

$ws.Range('D3').Value = '2.992.67'
$ws.Range('E3').Value = '  +2.08%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '381.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.96%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.07%  '

$ws.Range('E7').Value = '  +1.97%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.597'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.07%  '

$ws.Range('E10').Value = '  +1.96%  '

$ws.Range('E11').Value = '  -0.54%  '

$ws.Range('E12').Value = '  +1.89%  '

$ws.Range('D13').Value = '3.465.04'
$ws.Range('E13').Value = '  +2.24%  '

$ws.Range('E14').Value = '  +4.19%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '18.53'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.90%  '

$ws.Range('D16').Value = '2.990.81'
$ws.Range('E16').Value = '  +1.91%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '11.17'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.93%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.997'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.82%  '

$ws.Range('D19').Value = '51.625.65'
$ws.Range('E19').Value = '  +1.57%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.09'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.33%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.59'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.49%  '

$ws.Range('E22').Value = '  +1.30%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.50'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.28%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.81'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.11%  '

$ws.Range('E25').Value = '  +2.58%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.42%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.170'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.16%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.25'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.61%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.06%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '26.18'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.43%  '

$ws.Range('E31').Value = '  +1.02%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.40'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.78%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.64'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.62%  '

$ws.Range('E34').Value = '  +1.37%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.03'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.06%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0446'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.84%  '

$ws.Range('E37').Value = '  -0.06%  '

$ws.Range('E38').Value = '  +6.52%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.97'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.39%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.60'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.12%  '

$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '128.22'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.16%  '

$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.117'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.65%  '

$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.85'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.99%  '

$ws.Range('E44').Value = '  +14.07%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.41'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.83%  '

$ws.Range('E46').Value = '  -0.27%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.270'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.63%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.35'
$ws.Range('D48').Style = 'Normal'

$ws.Range('D49').Value = '2.037.39'
$ws.Range('E49').Value = '  +2.65%  '

$ws.Range('D50').Value = '3.289.84'
$ws.Range('E50').Value = '  +1.79%  '

$ws.Range('E51').Value = '  +1.81%  '
